$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

function Set-ParaText($textRange, $idx, $newText) {
    # Re-set to a short, distinct placeholder first. The text-diffing engine
    # preserves any common leading prefix (and its run formatting) between the
    # old and new paragraph text, which can split a single run into two when
    # old/new text happen to share a first character. Routing through a
    # placeholder that never shares a prefix with either string guarantees the
    # final assignment lands as one clean run with fresh rPr.
    $para = $textRange.Paragraphs($idx, 1)
    $para.Text = "|~placeholder~|"
    $para = $textRange.Paragraphs($idx, 1)
    $para.Text = $newText
}

Set-ParaText $tr 3 "Compare values and logically combine comparison results"
Set-ParaText $tr 4 "Manipulate lists of objects with a new, powerful datatype"
Set-ParaText $tr 5 "Conditionally execute or skip code"
Set-ParaText $tr 6 "Iterate over lists, strings and number sequences"
Set-ParaText $tr 7 "Loop over code until a condition is met"
Set-ParaText $tr 8 "Use some powerful built-in tools from the Python language to make iteration simple and expressive"

# Drop the old final bullet paragraph ("Start using the Python interpreter on
# your own.") entirely - its content got folded into paragraph 8 above.
$tr.Paragraphs(9, 1).Delete()

# The placeholder body shrinks once the (now shorter) bullet list no longer
# needs as much room; match the author's resulting box height exactly.
$sh.Height = 282.73443603515625
